# "adding lift to comparison"
#
# Adds a new "lift" data block (E20:E37) and a derived "norm_lift" column
# (F1:F18, = lift/125) to the lassos sheet, wires a 5th chart series
# ("norm_lift") into the existing line chart, repositions the chart /
# its legend, and restores the sheet selections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lassos")
$stats = $wb.Worksheets.Item("stats")

# ---------------------------------------------------------------------
# 1. Worksheet data: "lift" label + raw counts (E20:E37), "norm_lift"
#    header + formulas (F1:F18).
# ---------------------------------------------------------------------

# Shared-string insertion order matters for a faithful rebuild: "lift"
# is written first, then "norm_lift".
$ws.Range("E20").Value = "lift"

$liftValues = @(74, 72, 72, 69, 62, 61, 54, 51, 49, 48, 47, 42, 38, 36, 35, 34, 34)
for ($i = 0; $i -lt $liftValues.Length; $i++) {
    $ws.Cells.Item(21 + $i, 5).Value = $liftValues[$i]
}

$ws.Range("F1").Value = "norm_lift"
$ws.Range("F1").Font.Bold = $true

$ws.Range("F2").Formula = "=E21/125"
$ws.Range("F3:F18").Formula = "=E22/125"

# ---------------------------------------------------------------------
# 2. Chart: add the "norm_lift" series and flatten the smoothing that
#    the other series previously had.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
    $chart.SeriesCollection().Item($i).Smooth = $false
}

$newSeries = $chart.SeriesCollection().NewSeries()
$newSeries.Formula = "=SERIES(lassos!`$F`$1,lassos!`$A`$2:`$A`$18,lassos!`$F`$2:`$F`$18,5)"
$newSeries.Smooth = $false

# Legend moved from the top-right corner to the lower-left corner.
$legend = $chart.Legend
$legend.Left = 0.0583636103732116
$legend.Top = 0.643745992820316
$legend.Width = 0.234467818981024
$legend.Height = 0.323316984532656

# Chart frame grew / shifted right+down to make room for the new series.
$co.Left = 511.53515625
$co.Top = 55
$co.Width = 595.375
$co.Height = 533

# ---------------------------------------------------------------------
# 3. Restore view state (active cell selections on both sheets).
# ---------------------------------------------------------------------
$stats.Activate()
$stats.Range("D5").Select()

$ws.Activate()
$ws.Range("G33").Select()
